$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of row number -> new "Recorded By" (column G) value,
# reproducing the reordering of the comma-separated recorder list
# for the specific rows touched by the upstream sync.
$updates = @{
    2 = "backup@backdoor.com, system, System"
    4 = "backup@backdoor.com, System"
    5 = "backup@backdoor.com, System"
    8 = "backup@backdoor.com, System"
    11 = "System, dnasr281@gmail.com"
    17 = "System, dnasr281@gmail.com"
    28 = "backup@backdoor.com, system, System"
    30 = "backup@backdoor.com, System"
    31 = "backup@backdoor.com, System"
    34 = "backup@backdoor.com, System"
    37 = "System, dnasr281@gmail.com"
    43 = "System, dnasr281@gmail.com"
    54 = "backup@backdoor.com, system, System"
    56 = "backup@backdoor.com, System"
    57 = "backup@backdoor.com, System"
    60 = "backup@backdoor.com, System"
    63 = "System, dnasr281@gmail.com"
    69 = "System, dnasr281@gmail.com"
    80 = "backup@backdoor.com, System"
    81 = "backup@backdoor.com, System"
    82 = "backup@backdoor.com, System"
    93 = "System, dnasr281@gmail.com"
    94 = "System, dnasr281@gmail.com"
    96 = "System, dnasr281@gmail.com"
    106 = "backup@backdoor.com, System"
    107 = "backup@backdoor.com, System"
    108 = "backup@backdoor.com, System"
    119 = "System, dnasr281@gmail.com"
    120 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    132 = "backup@backdoor.com, System"
    133 = "backup@backdoor.com, System"
    134 = "backup@backdoor.com, System"
    145 = "System, dnasr281@gmail.com"
    146 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}

